$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header for column C (audioFalse -> currentPhase)
$ws.Range("C1").Value = "currentPhase"

# Update condition values in column C (rows 2 and 3) to the new shared value
$ws.Range("C2").Value = "train1P2"
$ws.Range("C3").Value = "train1P2"
